# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 170.08333
$ws.Range("I2").Value = 170.08333
$ws.Range("K2").Value = 170.08333
$ws.Range("M2").Value = -57.08332999999999
# Row 9
$ws.Range("H9").Value = 6569.2144
$ws.Range("I9").Value = 10171.111
$ws.Range("J9").Value = 85.8
$ws.Range("K9").Value = 10171.111
$ws.Range("L9").Value = 85.8
$ws.Range("M9").Value = -10002.111
$ws.Range("N9").Value = -423.8
# Row 12
$ws.Range("H12").Value = 135.47058
$ws.Range("I12").Value = 125.69231
$ws.Range("J12").Value = 167.25
$ws.Range("K12").Value = 125.69231
$ws.Range("L12").Value = 167.25
$ws.Range("M12").Value = 44.30768999999999
$ws.Range("N12").Value = -507.25
# Row 40
$ws.Range("H40").Value = 1832.6666
$ws.Range("I40").Value = 999
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 999
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -824
$ws.Range("N40").Value = -3850
# Row 98
$ws.Range("H98").Value = 1303.4073
$ws.Range("I98").Value = 991.6799999999999
$ws.Range("J98").Value = 5200
$ws.Range("K98").Value = 991.6799999999999
$ws.Range("L98").Value = 5200
$ws.Range("M98").Value = 506.3200000000001
$ws.Range("N98").Value = -8196
# Row 122
$ws.Range("H122").Value = 1303.4073
$ws.Range("I122").Value = 991.6799999999999
$ws.Range("J122").Value = 5200
$ws.Range("K122").Value = 2975.04
$ws.Range("L122").Value = 15600
$ws.Range("M122").Value = -525.04
$ws.Range("N122").Value = -20500
# Row 137
$ws.Range("H137").Value = 2109.2156
$ws.Range("J137").Value = 2094.2778
$ws.Range("L137").Value = 6282.8334
$ws.Range("N137").Value = -11382.8334
# Row 138
$ws.Range("H138").Value = 7146035.5
$ws.Range("J138").Value = 10207900
$ws.Range("L138").Value = 30623700
$ws.Range("N138").Value = -30633980
# Row 141
$ws.Range("H141").Value = 3518
$ws.Range("J141").Value = 3731.6667
$ws.Range("L141").Value = 11195.0001
$ws.Range("N141").Value = -21555.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 16166.667
$ws.Range("I37").Value = 18333.334
$ws.Range("J37").Value = 14000
$ws.Range("K37").Value = 18333.334
$ws.Range("L37").Value = 14000
$ws.Range("M37").Value = -18060.334
$ws.Range("N37").Value = -14546
# Row 74
$ws.Range("H74").Value = 12830.228
$ws.Range("I74").Value = 2633.6667
$ws.Range("J74").Value = 34680
$ws.Range("K74").Value = 2633.6667
$ws.Range("L74").Value = 34680
$ws.Range("M74").Value = -1759.6667
$ws.Range("N74").Value = -36428
# Row 77
$ws.Range("H77").Value = 12830.228
$ws.Range("I77").Value = 2633.6667
$ws.Range("J77").Value = 34680
$ws.Range("K77").Value = 13168.3335
$ws.Range("L77").Value = 173400
$ws.Range("M77").Value = -8800.333500000001
$ws.Range("N77").Value = -182136
# Row 122
$ws.Range("H122").Value = 1474.7333
$ws.Range("I122").Value = 1456.3334
$ws.Range("J122").Value = 1548.3334
$ws.Range("K122").Value = 4369.0002
$ws.Range("L122").Value = 4645.0002
$ws.Range("M122").Value = -1919.0002
$ws.Range("N122").Value = -9545.0002
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 132
$ws.Range("H132").Value = 3435.2144
$ws.Range("I132").Value = 3091.0833
$ws.Range("K132").Value = 9273.249899999999
$ws.Range("M132").Value = -6743.249899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3511.919
$ws.Range("I105").Value = 3585.8076
$ws.Range("J105").Value = 3337.2727
$ws.Range("K105").Value = 3585.8076
$ws.Range("L105").Value = 3337.2727
$ws.Range("M105").Value = -1838.8076
$ws.Range("N105").Value = -6831.2727
# Row 134
$ws.Range("H134").Value = 3323.9048
$ws.Range("I134").Value = 3116.9473
$ws.Range("K134").Value = 9350.841899999999
$ws.Range("M134").Value = -6815.841899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1950.5
$ws.Range("I31").Value = 1426.36
$ws.Range("J31").Value = 6318.3335
$ws.Range("K31").Value = 1426.36
$ws.Range("L31").Value = 6318.3335
$ws.Range("M31").Value = -1131.36
$ws.Range("N31").Value = -6908.3335
# Row 34
$ws.Range("H34").Value = 1950.5
$ws.Range("I34").Value = 1426.36
$ws.Range("J34").Value = 6318.3335
$ws.Range("K34").Value = 1426.36
$ws.Range("L34").Value = 6318.3335
$ws.Range("M34").Value = -1224.36
$ws.Range("N34").Value = -6722.3335
# Row 69
$ws.Range("H69").Value = 35000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 35000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -36498
# Row 72
$ws.Range("H72").Value = 35000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 105000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -112488
# Row 86
$ws.Range("H86").Value = 13858.95
$ws.Range("I86").Value = 14893.625
$ws.Range("K86").Value = 14893.625
$ws.Range("M86").Value = -13770.625
# Row 87
$ws.Range("H87").Value = 38750
$ws.Range("J87").Value = 38750
$ws.Range("L87").Value = 38750
$ws.Range("N87").Value = -41122
# Row 89
$ws.Range("H89").Value = 13858.95
$ws.Range("I89").Value = 14893.625
$ws.Range("K89").Value = 74468.125
$ws.Range("M89").Value = -68852.125
# Row 90
$ws.Range("H90").Value = 38750
$ws.Range("J90").Value = 38750
$ws.Range("L90").Value = 116250
$ws.Range("N90").Value = -128106
# Row 141
$ws.Range("H141").Value = 219722.86
$ws.Range("J141").Value = 219722.86
$ws.Range("L141").Value = 219722.86
$ws.Range("N141").Value = -230082.86

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 19.3
$ws.Range("I2").Value = 16.625
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 99.75
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = 13.25
$ws.Range("N2").Value = -406

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 43300000
$ws.Range("I10").Value = 53500000
$ws.Range("J10").Value = 2500000
$ws.Range("K10").Value = 53500000
$ws.Range("L10").Value = 2500000
$ws.Range("M10").Value = -53499831
$ws.Range("N10").Value = -2500338
# Row 80
$ws.Range("H80").Value = 3595.8
$ws.Range("I80").Value = 3660
$ws.Range("J80").Value = 3499.5
$ws.Range("K80").Value = 3660
$ws.Range("L80").Value = 3499.5
$ws.Range("M80").Value = -2662
$ws.Range("N80").Value = -5495.5
# Row 83
$ws.Range("H83").Value = 3595.8
$ws.Range("I83").Value = 3660
$ws.Range("J83").Value = 3499.5
$ws.Range("K83").Value = 18300
$ws.Range("L83").Value = 17497.5
$ws.Range("M83").Value = -13308
$ws.Range("N83").Value = -27481.5
# Row 102
$ws.Range("H102").Value = 45455680
$ws.Range("I102").Value = 1180.5
$ws.Range("K102").Value = 1180.5
$ws.Range("M102").Value = 441.5
# Row 113
$ws.Range("H113").Value = 4148.2144
$ws.Range("I113").Value = 4097.9165
$ws.Range("K113").Value = 4097.9165
$ws.Range("M113").Value = -1927.9165
# Row 132
$ws.Range("H132").Value = 4369.931
$ws.Range("J132").Value = 5483.5
$ws.Range("L132").Value = 16450.5
$ws.Range("N132").Value = -21510.5
# Row 139
$ws.Range("H139").Value = 116324
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4637.773
$ws.Range("I40").Value = 4218.7334
$ws.Range("J40").Value = 5535.7144
$ws.Range("K40").Value = 4218.7334
$ws.Range("L40").Value = 5535.7144
$ws.Range("M40").Value = -4082.7334
$ws.Range("N40").Value = -5807.7144
# Row 61
$ws.Range("H61").Value = 4940.2
$ws.Range("I61").Value = 4417
$ws.Range("J61").Value = 5725
$ws.Range("K61").Value = 4417
$ws.Range("L61").Value = 5725
$ws.Range("M61").Value = -4215
$ws.Range("N61").Value = -6129
# Row 93
$ws.Range("H93").Value = 1973.6296
$ws.Range("J93").Value = 1786
$ws.Range("L93").Value = 1786
$ws.Range("N93").Value = -4282
# Row 113
$ws.Range("H113").Value = 4940.2
$ws.Range("I113").Value = 4417
$ws.Range("J113").Value = 5725
$ws.Range("K113").Value = 4417
$ws.Range("L113").Value = 5725
$ws.Range("M113").Value = -2247
$ws.Range("N113").Value = -10065
# Row 132
$ws.Range("H132").Value = 104680.2
$ws.Range("I132").Value = 129225
$ws.Range("J132").Value = 6501
$ws.Range("K132").Value = 387675
$ws.Range("L132").Value = 19503
$ws.Range("M132").Value = -385145
$ws.Range("N132").Value = -24563
# Row 136
$ws.Range("H136").Value = 5087.375
$ws.Range("I136").Value = 4999.909
$ws.Range("J136").Value = 5279.8
$ws.Range("K136").Value = 14999.727
$ws.Range("L136").Value = 15839.4
$ws.Range("M136").Value = -12449.727
$ws.Range("N136").Value = -20939.4

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 126999.875
$ws.Range("I107").Value = 1999.8
$ws.Range("K107").Value = 5999.4
$ws.Range("M107").Value = -4079.4
# Row 113
$ws.Range("H113").Value = 1354
$ws.Range("I113").Value = 1807.75
$ws.Range("K113").Value = 5423.25
$ws.Range("M113").Value = -3253.25
# Row 122
$ws.Range("H122").Value = 2670.2693
$ws.Range("I122").Value = 2517.7917
$ws.Range("K122").Value = 7553.375100000001
$ws.Range("M122").Value = -5103.375100000001
# Row 131
$ws.Range("H131").Value = 154994
$ws.Range("J131").Value = 154994
$ws.Range("L131").Value = 154994
$ws.Range("N131").Value = -165074
# Row 132
$ws.Range("H132").Value = 3272.2727
$ws.Range("I132").Value = 2887.7778
$ws.Range("J132").Value = 5002.5
$ws.Range("K132").Value = 8663.3334
$ws.Range("L132").Value = 15007.5
$ws.Range("M132").Value = -6133.3334
$ws.Range("N132").Value = -20067.5
